# The sheet lists people with a "CvLAC" profile URL in column E.
# Row 2 (Adolfo Alarcon Guzman) had the wrong CvLAC link (cod_rh=0000005410,
# which actually belongs to a different/bad record - the commit message refers
# to cleaning up a bad char(") on names, i.e. fixing a previously-mismatched
# row). Re-point E2 to the correct profile URL (cod_rh=0000218430) and make it
# a real hyperlink, the way Excel does when you paste/insert a link by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctUrl = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000218430"

$cell = $ws.Range("E2")

# Adding the hyperlink also applies the built-in "Hyperlink" cell style
# (blue, underlined) to the cell - exactly like Excel's Insert Hyperlink.
$ws.Hyperlinks.Add($cell, $correctUrl)

# Make sure the cell's displayed text is the new URL as well.
$cell.Value = $correctUrl

# Reflect the new selection left after the edit (Excel leaves the
# just-edited cell selected).
$cell.Select()
